$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("quiz")

# Update "Marking" row - Right column
$ws.Range("B11").Value = 5

# Update "Total" row - Right column and Max (correct/total marks) column
$ws.Range("B12").Value = 90
$ws.Range("E12").Value = "90/140"
